$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-11) holds a date-formatted serial number that needs to
# move forward by one day: 45205 (2023-10-06) -> 45206 (2023-10-07).
for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2()
    if ([string]$current -eq "45205") {
        $cell.Value = 45206
    }
}
